$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current last data row (old row 17), pushing
# the old row 17 down to row 19 and the closing-signature rows (22,23) down
# to rows 24,25 - matching the target layout.
$ws.Rows("17:18").Insert()

# Populate the two newly inserted rows with the same formatting as row 16
# (the existing "Salario Basico" detail row) by copying it down twice; we
# overwrite the actual values afterwards.
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))

# Row 18 keeps exactly what row 16 used to contain (YUDIE PADILLA CARRASQUILLA),
# so it is already correct after the copy above.

# Row 16 now becomes the new worker record (NILKA PAOLA GOMEZ PEREZ, period 2005)
$ws.Range("C16").Value2 = "45529300"
$ws.Range("D16").Value2 = "NILKA PAOLA GOMEZ PEREZ"
$ws.Range("E16").Value2 = "2005"
$ws.Range("F16").Value2 = 35112
$ws.Range("G16").Value2 = 877803

# Row 17 repeats the same new worker for a second period (2004)
$ws.Range("C17").Value2 = "45529300"
$ws.Range("D17").Value2 = "NILKA PAOLA GOMEZ PEREZ"
$ws.Range("E17").Value2 = "2004"
$ws.Range("F17").Value2 = 35112
$ws.Range("G17").Value2 = 877803

# Update the account summary figures at the top of the sheet
$ws.Range("E11").Value2 = 116009
$ws.Range("C13").Value2 = 3
$ws.Range("F13").Value2 = 4
